$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 27: num_customers 46 -> 47, retention_rate recalculated (cohort_size stays 2252)
$ws.Range("C27").Value = 47
$ws.Range("E27").Value = 0.02087033747779751

# Row 37: num_customers 822 -> 825, cohort_size 822 -> 825 (retention_rate stays 1)
$ws.Range("C37").Value = 825
$ws.Range("D37").Value = 825
